# Update the "Attendance Roster" workbook: Sheet2 gains a third IP column (IP3)
# and the IP1/IP2/IP3 address lists are refreshed with a new batch of
# AWS host addresses (old batch removed from sharedStrings, new batch added).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Header row (row 11): IP1 / IP2 / IP3 / Username -------------------------
$ws.Range("G11").Value = "IP1"
$ws.Range("H11").Value = "IP2"
# I11 used to hold the "Username" header; it now becomes the new IP3 header.
# Give it the same look as the other header cells by cloning the format from
# the existing header cell before writing the new text.
$ws.Range("G11").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = "IP3"
$ws.Range("J11").PasteSpecial(-4122)
$ws.Range("J11").Value = "Username"
$excel.CutCopyMode = 0

# --- Data rows 12-23: new IP1 / IP2 / IP3 / Username values ------------------
$rows = @(
    @{ Row = 12; G = "100.27.231.78  "; H = "23.20.137.133";  I = "54.91.79.5" },
    @{ Row = 13; G = "54.167.4.16";     H = "54.90.202.68";   I = "54.159.68.216" },
    @{ Row = 14; G = "98.91.18.148";    H = "18.234.65.11";   I = "3.80.186.244" },
    @{ Row = 15; G = "44.222.199.240";  H = "34.228.195.129"; I = "54.221.66.63" },
    @{ Row = 16; G = "54.242.87.139";   H = "34.227.225.224"; I = "184.72.209.116" },
    @{ Row = 17; G = "23.23.38.212";    H = "54.209.249.58";  I = "98.84.174.181" },
    @{ Row = 18; G = "52.23.158.12";    H = "13.221.126.39";  I = "54.82.99.218" },
    @{ Row = 19; G = "13.222.175.186";  H = "54.87.32.229";   I = "3.80.32.155" },
    @{ Row = 20; G = "54.210.132.16";   H = "54.242.46.158";  I = "3.82.52.45" },
    @{ Row = 21; G = "54.87.220.154";   H = "3.95.220.157";   I = "184.72.209.49" },
    @{ Row = 22; G = "34.224.40.166";   H = "3.82.37.37";     I = "54.82.67.219" },
    @{ Row = 23; G = "23.22.133.70";    H = "34.229.90.9";    I = "54.242.131.20" }
)

foreach ($r in $rows) {
    $ws.Range("G$($r.Row)").Value = $r.G
    $ws.Range("H$($r.Row)").Value = $r.H
    $ws.Range("I$($r.Row)").Value = $r.I
}

# New column J holds the "ubuntu" username for every data row, formatted the
# same way as the adjoining I column.
$ws.Range("I12:I23").Copy()
$ws.Range("J12:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J12:J23").Value = "ubuntu"

# --- Column widths -------------------------------------------------------
# (columns are addressed by number - Columns.Item("G") does not resolve)
$ws.Columns.Item(7).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 13
$ws.Columns.Item(9).ColumnWidth = 13
$ws.Columns.Item(10).ColumnWidth = 9
$ws.Columns.Item(11).ColumnWidth = 14.166666666666666

# --- Selection -------------------------------------------------------------
$ws.Range("K16").Select()
